$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: advance the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value2 = $ws.Range("A1").Value2 + 1

# Step 2: update the two price cells
$ws.Range("D28").Value = 42870
$ws.Range("D29").Value = 57191
